$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02638114373498337
$ws.Range("C2").Value = 0.2136373979745126
$ws.Range("D2").Value = 0.05862372610272825
$ws.Range("E2").Value = 0.2421233695922974
$ws.Range("F2").Value = 0.2497673777755665

$ws.Range("B3").Value = 0.05401209062359413
$ws.Range("C3").Value = 0.218968316772212
$ws.Range("D3").Value = 0.0812101156932476
$ws.Range("E3").Value = 0.284973885984747
$ws.Range("F3").Value = 0.2912339447472185

$ws.Range("B4").Value = 0.04189333505767585
$ws.Range("C4").Value = 0.253872083886316
$ws.Range("D4").Value = 0.1125880168289207
$ws.Range("E4").Value = 0.3355413787134467
$ws.Range("F4").Value = 0.3477192681250979

$ws.Range("B5").Value = -0.0001372508960761152
$ws.Range("C5").Value = 0.2480625890636292
$ws.Range("D5").Value = 0.08457529259376213
$ws.Range("E5").Value = 0.2908183154372539
$ws.Range("F5").Value = 0.3050127884721378

$ws.Range("B6").Value = -0.01622121920471854
$ws.Range("C6").Value = 0.2631740310991668
$ws.Range("D6").Value = 0.09576598071011208
$ws.Range("E6").Value = 0.3094607902628572
$ws.Range("F6").Value = 0.325751870051125

$ws.Range("B7").Value = -0.03331298675845948
$ws.Range("C7").Value = 0.2786008952859811
$ws.Range("D7").Value = 0.1208052536752846
$ws.Range("E7").Value = 0.3475705017335111
$ws.Range("F7").Value = 0.3669569946357198

$ws.Range("B8").Value = -0.123504592765017
$ws.Range("C8").Value = 0.2639331087351315
$ws.Range("D8").Value = 0.08609388686166163
$ws.Range("E8").Value = 0.2934175980776573
$ws.Range("F8").Value = 0.2915623482432715

$ws.Range("B9").Value = -0.2310012636481934
$ws.Range("C9").Value = 0.2891182804648226
$ws.Range("D9").Value = 0.1126651299999655
$ws.Range("E9").Value = 0.3356562676309881
$ws.Range("F9").Value = 0.2982537833613431

$ws.Range("B10").Value = -0.1742297805489477
$ws.Range("C10").Value = 0.1742297805489477
$ws.Range("D10").Value = 0.03035601643013447
$ws.Range("E10").Value = 0.1742297805489477
